$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "sum" header (G1) onto the new H1
# header cell so the new "Save" header matches the bold/centered/bordered
# header style used by the rest of row 1, then set its text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New data column value for row 2.
$ws.Range("H2").Value = 0
